# Generate Report for Handoff
# Status moves from "In Translation" to "Ready for handoff" and the
# handoff timestamps advance by 30s. The "Status" columns are widened
# (they auto-size to the new, longer status text).

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-07 17:26:49"
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3

# ---- zh-cn sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-07 17:26:44"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3

# ---- de-de sheet ---------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-07 17:26:49"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3

Write-Host "Applied handoff-ready status update"
